# Update the "Förändrad" (changed) date column (C2:C5) from 2023-09-15 to 2023-09-16
# Corresponds to serial date value change 45184 -> 45185

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = Get-Date -Year 2023 -Month 9 -Day 16 -Hour 0 -Minute 0 -Second 0

$ws.Range("C2").Value = $newDate
$ws.Range("C3").Value = $newDate
$ws.Range("C4").Value = $newDate
$ws.Range("C5").Value = $newDate
